# 20240919 Completado Descargar Excel, Descargar Gráficos en EBA y Metricas
#
# The "Cost to income ratio" row (row 11) had its Spanish translation
# renamed from "Ratio de costos a ingresos" to "Ratio de eficiencia".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("listado_nombres_Yahoo")

$ws.Range("B11").Value = "Ratio de eficiencia"

# Restore the active selection to where the user left off editing.
$ws.Activate()
$ws.Range("B12").Select() | Out-Null
